# "Added Error Comments in RowMode"
# Adds a new Sheet2 after Sheet1, mirroring the layout of Sheet1's A1:D5
# block (transposed to M14:P18) but with the numeric "error comment" style
# values (0) in column N and a distinguishing big number (50000) in O17.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create Sheet2 right after Sheet1 -----------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Header row (row 14) -------------------------------------------------
$ws2.Range("M14").Value = "StringTest"
$ws2.Range("N14").Value = "DecimalTest"
$ws2.Range("O14").Value = "IntTest"
$ws2.Range("P14").Value = "GuidTest"

# --- Row 15 ---------------------------------------------------------------
$ws2.Range("M15").Value = "A"
$ws2.Range("N15").Value = 0
$ws2.Range("O15").Value = 1
$ws2.Range("P15").Value = "291E645A-F0A6-4A50-A316-7A7A9931C8F2"
$ws2.Range("P15").NumberFormat = "0.00E+00"

# --- Row 16 ---------------------------------------------------------------
$ws2.Range("M16").Value = "B"
$ws2.Range("N16").Value = 0
$ws2.Range("O16").Value = 4
$ws2.Range("P16").Value = "9BAB1466-A66A-4DB7-9AD4-DA1ED9531193"

# --- Row 17 ---------------------------------------------------------------
$ws2.Range("M17").Value = "Ä"
$ws2.Range("N17").Value = "1,5"
$ws2.Range("N17").NumberFormat = "@"
$ws2.Range("O17").Value = 50000
$ws2.Range("P17").Value = "6B10C320-B126-43F3-904E-17571550AD16"

# --- Row 18 ---------------------------------------------------------------
$ws2.Range("M18").Value = "D"
$ws2.Range("N18").Value = 0
$ws2.Range("O18").Value = 7
$ws2.Range("P18").Value = "2F93BE1B-433F-41FE-9B3E-22452F18F20C"

# --- View/selection state --------------------------------------------------
# Sheet1 keeps a selection over its data block (A1:D5, active cell D5).
$ws1.Activate()
$ws1.Range("A1:D5").Select()

# Sheet2 becomes the active tab, with N17 as the active cell.
$ws2.Activate()
$ws2.Range("N17").Select()
